$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 240, shifting existing rows 240:320 down to 241:321
$ws.Rows.Item(240).Insert()

# Populate the newly inserted row 240 with the new record
$ws.Cells.Item(240, 1).Value = 5
$ws.Cells.Item(240, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(240, 3).Value = "Maule"
$ws.Cells.Item(240, 4).Value = 44985
$ws.Cells.Item(240, 5).Value = 7
$ws.Cells.Item(240, 6).Value = 100112024
$ws.Cells.Item(240, 7).Value = "Choclo"
$ws.Cells.Item(240, 8).Value = "Choclero"
$ws.Cells.Item(240, 9).Value = "Primera"
$ws.Cells.Item(240, 10).Value = 20000
$ws.Cells.Item(240, 11).Value = 500
$ws.Cells.Item(240, 12).Value = 500
$ws.Cells.Item(240, 13).Value = 500
$ws.Cells.Item(240, 14).Value = "`$/unidad"
$ws.Cells.Item(240, 15).Value = "Región del Maule"
$ws.Cells.Item(240, 16).Value = 500
$ws.Cells.Item(240, 17).Value = 1
$ws.Cells.Item(240, 18).Value = "Hortaliza"
